# Fix arrays not being filled
# Changed how arrays are constructed / filled arrays with data
#
# The " @stringArray" settings row (row 10, cols E:I) together with its
# sample values "string_0".." string_4" (row 11, cols E:I) were leftover
# placeholders from the old (broken) array construction. They are
# replaced with a proper " @intArray" setting and its filled-in
# "int_0".."int_4" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: label for the array setting (was " @stringArray" repeated, now " @intArray")
$ws.Range("E10").Value = " @intArray"
$ws.Range("F10").Value = " @intArray"
$ws.Range("G10").Value = " @intArray"
$ws.Range("H10").Value = " @intArray"
$ws.Range("I10").Value = " @intArray"

# Row 11: the actual array values that were missing/wrong (string_0..string_4 -> int_0..int_4)
$ws.Range("E11").Value = "int_0"
$ws.Range("F11").Value = "int_1"
$ws.Range("G11").Value = "int_2"
$ws.Range("H11").Value = "int_3"
$ws.Range("I11").Value = "int_4"

# Update the view state: scroll down so row 7 is at the top and select I21
# (matches where the user ended up after filling in the array).
$ws.Range("I21").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
